$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# --- Crime-stat table updates (rows 15-31) ---
$ws.Range("F15").Copy($ws.Range("C15"))
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 75
$ws.Range("L15").Value = 55.555555555555
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 600
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 183.333333333333
$ws.Range("I16").Value = 77
$ws.Range("J16").Value = 68
$ws.Range("K16").Value = 13.235294117647
$ws.Range("L16").Value = 35.087719298245
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -48
$ws.Range("I17").Value = 159
$ws.Range("J17").Value = 186
$ws.Range("K17").Value = -14.516129032258
$ws.Range("L17").Value = -7.558139534883
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 75
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = 50
$ws.Range("L18").Value = 11.940298507462
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.76923076923
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 278
$ws.Range("J19").Value = 328
$ws.Range("K19").Value = -15.243902439024
$ws.Range("L19").Value = -12.302839116719
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 72
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = -15.294117647058
$ws.Range("L20").Value = -6.493506493506
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = 3.614457831325
$ws.Range("I21").Value = 677
$ws.Range("J21").Value = 729
$ws.Range("K21").Value = -7.13305898491
$ws.Range("L21").Value = -3.147353361945
$ws.Range("D23").Copy($ws.Range("C23"))
$ws.Range("F23").Value = 1
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = -15.447154471544
$ws.Range("I24").Value = 1090
$ws.Range("J24").Value = 1108
$ws.Range("K24").Value = -1.624548736462
$ws.Range("L24").Value = 0
$ws.Range("C25").Value = 20
$ws.Range("E25").Value = -9.090909090909
$ws.Range("F25").Value = 72
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 4.347826086956
$ws.Range("I25").Value = 737
$ws.Range("J25").Value = 649
$ws.Range("K25").Value = 13.559322033898
$ws.Range("L25").Value = 40.917782026768
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 14.285714285714
$ws.Range("I26").Value = 454
$ws.Range("J26").Value = 399
$ws.Range("K26").Value = 13.784461152882
$ws.Range("L26").Value = 20.10582010582
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 60
$ws.Range("F28").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 5
$ws.Range("D29").Copy($ws.Range("D28"))
$ws.Range("E29").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 49
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 36.111111111111
$ws.Range("F31").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("H31").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 5
$ws.Range("K31").Value = 40
